# Update gh-pages output (杭州-漫展信息.xlsx) as of commit 456a3b4
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Sheet 1: 展览 (Exhibitions)
# ---------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)

$ws1.Range("F2").Value = 1293
$ws1.Range("G4").Value = "不可售"
$ws1.Range("F6").Value = 6886
$ws1.Range("F7").Value = 1830
$ws1.Range("F8").Value = 6419
$ws1.Range("F9").Value = 148
$ws1.Range("F10").Value = 1974
$ws1.Range("F11").Value = 525
$ws1.Range("F12").Value = 24
$ws1.Range("F17").Value = 57
$ws1.Range("F18").Value = 8143
$ws1.Range("F19").Value = 145
$ws1.Range("F22").Value = 111
$ws1.Range("F23").Value = 1767
$ws1.Range("F24").Value = 851
$ws1.Range("F28").Value = 49
$ws1.Range("F29").Value = 178
$ws1.Range("F31").Value = 1874
$ws1.Range("F32").Value = 822
$ws1.Range("F33").Value = 405
$ws1.Range("F36").Value = 129
$ws1.Range("F37").Value = 93
$ws1.Range("F39").Value = 3937

# ---------------------------------------------------------------
# Sheet 2: 演出 (Performances)
# ---------------------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)

$ws2.Range("F4").Value = 376

# ---------------------------------------------------------------
# Sheet 3: 本地生活 (Local life)
# ---------------------------------------------------------------
$ws3 = $wb.Worksheets.Item(3)

$ws3.Range("F3").Value = 2298
$ws3.Range("F4").Value = 696

# ---------------------------------------------------------------
# Sheet 4: 全部类型 (All types) - combined view
# ---------------------------------------------------------------
$ws4 = $wb.Worksheets.Item(4)

$ws4.Range("F3").Value = 2298
$ws4.Range("F4").Value = 696
$ws4.Range("F5").Value = 1293

# Row 7 becomes the former row-8 event (2024CJMF), content shifted down / replaced
# (force text so Excel doesn't auto-convert the date-like string to a serial date)
$ws4.Range("B7").NumberFormat = "@"
$ws4.Range("B7").Value = "2024-09-15"
$ws4.Range("B7").ClearFormats()
$ws4.Range("C7").Value = "杭州·2024CJMF·不止音乐节"
$ws4.Range("D7").Value = "塘子堰路177号 华数产业园隔壁大草坪"
$ws4.Range("E7").Value = "2024.09.15 13:00-09.16 21:40"
$ws4.Range("F7").Value = 376
$ws4.Range("G7").Value = 168
$ws4.Range("H7").Value = "https://show.bilibili.com/platform/detail.html?id=90522"
$ws4.Range("I7").Value = "//i1.hdslb.com/bfs/openplatform/202408/3PmG2Bq51723192884141.jpeg"

# Row 8 becomes the former row-9 event (COMIC GALAXY)
$ws4.Range("C8").Value = "杭州·2024首届COMIC GALAXY次元盛典"
$ws4.Range("D8").Value = "长江南路336号 白马湖国际会展中心"
$ws4.Range("E8").Value = "2024.09.15 09:30-09.17 17:30"
$ws4.Range("F8").Value = 6886
$ws4.Range("G8").Value = 88
$ws4.Range("H8").Value = "https://show.bilibili.com/platform/detail.html?id=90433"
$ws4.Range("I8").Value = "//i0.hdslb.com/bfs/openplatform/202408/teoBMbzd1723019674766.png"

# Row 9 becomes the former row-10 event (伤心咖啡馆之歌 stage play)
$ws4.Range("C9").Value = "杭州·多感官环境式话剧《伤心咖啡馆之歌》"
$ws4.Range("D9").Value = "留泗路东山里22号 大美创意园-2号楼"
$ws4.Range("E9").Value = "2024.09.15 15:00-09.17 21:00"
$ws4.Range("F9").Value = 1
$ws4.Range("G9").Value = 78
$ws4.Range("H9").Value = "https://show.bilibili.com/platform/detail.html?id=91995"
$ws4.Range("I9").Value = "//i0.hdslb.com/bfs/openplatform/202409/9PRTu6Fm1725437709663.png"

# Row 10 becomes a new event (木灵动漫 二哈和他的白猫师尊主题餐厅)
$ws4.Range("C10").Value = "杭州·木灵动漫 二哈和他的白猫师尊主题餐厅"
$ws4.Range("D10").Value = "平海路124号 杭州湖滨88"
$ws4.Range("E10").Value = "2024.09.15 00:00-09.30 23:59"
$ws4.Range("F10").Value = 286
$ws4.Range("G10").Value = 10
$ws4.Range("H10").Value = "https://show.bilibili.com/platform/detail.html?id=91251"
$ws4.Range("I10").Value = "//i2.hdslb.com/bfs/openplatform/202408/wLlo6EFv1724642759732.png"

$ws4.Range("F11").Value = 1830
$ws4.Range("F12").Value = 6419
$ws4.Range("F13").Value = 148
$ws4.Range("F14").Value = 1974
$ws4.Range("F16").Value = 525
$ws4.Range("F23").Value = 57
$ws4.Range("F24").Value = 8143
$ws4.Range("F25").Value = 145
$ws4.Range("F28").Value = 111
$ws4.Range("F29").Value = 1767
$ws4.Range("F30").Value = 851
$ws4.Range("F33").Value = 178
$ws4.Range("F34").Value = 1874
$ws4.Range("F35").Value = 822
$ws4.Range("F37").Value = 405
$ws4.Range("F44").Value = 3937

$wb.Save()
